$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,1,1,0,3,1,3,0,1,1,2,2,0,2,0,1,5,0,1,2,3,6,1,6,1,3,5,1,1,1,1,1,2,5,4,1,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
